# Updated cryptos list on Sun Oct  8 14:18:47 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) values for
# each coin row (rows 2-50) to the latest scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D/E to text format so that numeric-looking
# strings (e.g. "210.86", "23.31") are written back as plain text,
# matching the original inline-string cell contents instead of being
# auto-converted into floating point numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.857.52"
$ws.Range("D3").Value = "1.625.20"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "210.86"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D8").Value = "23.31"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "1.855.83"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "1.626.29"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "0.561"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").Value = "65.21"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "27.852.03"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "229.32"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "7.65"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "10.09"
$ws.Range("E23").Value = "  -3.26%  "
$ws.Range("D24").Value = "2.06"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").Value = "154.16"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "6.88"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "15.50"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").Value = "3.41"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").Value = "1.395.23"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D36").Value = "1.03"
$ws.Range("E36").Value = "  +11.37%  "
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "0.555"
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "65.69"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("D45").Value = "5.43"
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("D46").Value = "1.767.57"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("E47").Value = "  -3.09%  "
$ws.Range("D48").Value = "87.94"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").Value = "0.0₆0101"
$ws.Range("E50").Value = "  -3.89%  "

# Restore the cells original (default) style now that the values are set.
$dataRange.Style = "Normal"
